# Weekly Fruta/Hortaliza update:
# - Re-orders the existing weekly rows (2-12) to their new values per the
#   latest data pull (dates + volumen/precio/precio-kg columns).
# - Appends a new weekly row (13) for the "Granada" record that is now
#   the most recent one in the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 2-12: Fecha(D), Volumen(M),
# Precio minimo(N), Precio maximo(O), Precio promedio ponderado(P),
# Precio $/Kg(S). All other columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) are
# unchanged for these rows.
$updates = @(
    @{ Row = 2;  D = 45096; M = 30; N = 20000; O = 20000; P = 20000; S = 1111 },
    @{ Row = 3;  D = 45083; M = 50; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 4;  D = 45069; M = 60; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 5;  D = 45084; M = 50; N = 18000; O = 19000; P = 18500; S = 1028 },
    @{ Row = 6;  D = 45076; M = 20; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 7;  D = 45061; M = 40; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 9;  D = 45085; M = 30; N = 19000; O = 19000; P = 19000; S = 1056 },
    @{ Row = 10; D = 45055; M = 50; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 11; D = 45111; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 },
    @{ Row = 12; D = 45112; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 13).Value = $u.M
    $ws.Cells.Item($r, 14).Value = $u.N
    $ws.Cells.Item($r, 15).Value = $u.O
    $ws.Cells.Item($r, 16).Value = $u.P
    $ws.Cells.Item($r, 19).Value = $u.S
}

# Row 13 is a brand-new record; copy the formatting of row 12 (so the date
# cell keeps the same date/time number format) then fill in all values.
$ws.Range("A12:T12").Copy()
$ws.Range("A13:T13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 45072
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100104
$ws.Cells.Item(13, 8).Value = "Frutos de pepita"
$ws.Cells.Item(13, 9).Value = 100104001
$ws.Cells.Item(13, 10).Value = "Granada"
$ws.Cells.Item(13, 11).Value = "Wonderfull"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 30
$ws.Cells.Item(13, 14).Value = 15000
$ws.Cells.Item(13, 15).Value = 15000
$ws.Cells.Item(13, 16).Value = 15000
$ws.Cells.Item(13, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(13, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(13, 19).Value = 833
$ws.Cells.Item(13, 20).Value = 18
